# SiouxFallNet/processdata.xlsx
# Commit: "create two sioux data"
#
# Adds a new Sheet4 that mirrors Sheet1's A/B node-pair data (76 rows) and
# augments it with two helper formula columns (C = A+1, D = B+1), highlights
# one row pair (A38:D39) with a yellow fill, and updates the view state so
# Sheet4 becomes the active/selected sheet (Sheet3 loses its old selection).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new worksheet after the last existing sheet (Sheet3) so it
#    becomes "Sheet4".
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)

# ---------------------------------------------------------------------------
# 2. Column A / B data -- identical node-pair list to Sheet1!A1:B76.
# ---------------------------------------------------------------------------
$aVals = @(0,0,1,1,2,2,2,3,3,3,4,4,4,5,5,5,6,6,7,7,7,7,8,8,8,9,9,9,9,9,10,10,10,10,11,11,11,12,12,13,13,13,14,14,14,14,15,15,15,15,16,16,16,17,17,17,18,18,18,19,19,19,19,20,20,20,21,21,21,21,22,22,22,23,23,23)
$bVals = @(1,2,0,5,0,3,11,2,4,10,3,5,8,1,4,7,7,17,5,6,8,15,4,7,9,8,10,14,15,16,3,9,11,13,2,10,12,11,23,10,14,22,9,13,18,21,7,9,16,17,9,15,18,6,15,19,14,16,19,17,18,20,21,19,21,23,14,19,20,22,13,21,23,12,20,22)

$rowCount = $aVals.Count
$data = New-Object 'object[,]' $rowCount,2
for ($i = 0; $i -lt $rowCount; $i++) {
    $data[$i,0] = $aVals[$i]
    $data[$i,1] = $bVals[$i]
}
$ws4.Range("A1:B$rowCount").Value = $data

# ---------------------------------------------------------------------------
# 3. Column C / D -- "+1" helper formulas for every data row.
# ---------------------------------------------------------------------------
$ws4.Range("C1:C$rowCount").Formula = "=A1+1"
$ws4.Range("D1:D$rowCount").Formula = "=B1+1"

# ---------------------------------------------------------------------------
# 4. Highlight A38:D39 with a solid yellow fill.
# ---------------------------------------------------------------------------
$ws4.Range("A38:D39").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 5. View state: Sheet3 loses its old selection/active status, Sheet4 gets
#    selected/scrolled to the highlighted rows and becomes the active sheet.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
[void]$ws3.Range("B550").Select()

$ws4.Activate()
[void]$ws4.Range("A38:B39").Select()
$excel.ActiveWindow.ScrollRow = 28
